# edit.ps1
# Applies three changes to MatteoRiondato-CV.docx:
#  1. Add a new "Journal Articles" entry for the ABRA TKDD journal paper
#     (accepted, to appear), right before the TRIEST/TKDD 2017 entry.
#  2. Add "'18" to the "IEEE ICDM" Program Committees line.
#  3. Update the cached footer PAGE field result from "7" to "2".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the new "Accepted ... ABRA ... To appear" journal paragraph
#    immediately before the paragraph that begins "2017" and discusses
#    "Counting Local and Global Triangles" (the TRIEST / TKDD entry).
# ---------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "2017*" -and $t -like "*Counting Local and Global Triangles*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate the TRIEST/TKDD 2017 journal paragraph"
}

$p = $d.Paragraphs.Item($targetIndex)
$rng = $p.Range
$rng.Collapse(1)
$frag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="normal0"/><w:spacing w:after="120"/><w:ind w:left="1440" w:hanging="1440"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t>Accepted</w:t></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">M. Riondato </w:t></w:r><w:r><w:t>and E. Upfal. ABRA: Approximating Betweenness Centrality in Static and Dynamic Graphs with Rademacher Averages</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>ACM Transactions on Knowledge Discovery from Data</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> To appear</w:t></w:r></w:p><w:p w:rsidR="00F9559A" w:rsidRPr="00252DBF" w:rsidRDefault="00252DBF"><w:pPr><w:pStyle w:val="normal0"/><w:spacing w:after="120"/><w:ind w:left="1440" w:hanging="1440"/></w:pPr><w:r><w:t>2017</w:t></w:r><w:r w:rsidR="00F9559A"><w:tab/></w:r><w:r w:rsidR="00F9559A" w:rsidRPr="00E77A52"><w:t xml:space="preserve">L. De Stefani, A. Epasto, </w:t></w:r><w:r w:rsidR="00F9559A" w:rsidRPr="00E77A52"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>M. Riondato</w:t></w:r><w:r w:rsidR="00F9559A" w:rsidRPr="00E77A52"><w:t>, and E. Upfal. TRIÉST: Counting Local and Global Triangles in Fully-dynamic Streams with Fixed Memory Size.</w:t></w:r><w:r w:rsidR="00F9559A"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00F9559A"><w:rPr><w:i/></w:rPr><w:t>ACM Transactions on Knowledge Discovery from Data</w:t></w:r><w:r><w:t>, 11(4):43:1</w:t></w:r><w:r w:rsidRPr="00E77A52"><w:t>–</w:t></w:r><w:r><w:t>43:50</w:t></w:r><w:r w:rsidR="005C454F"><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="005C454F"><w:rPr><w:b/></w:rPr><w:t>Invited article to the special issue on the best papers from KDD 2016</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($frag1)

# ---------------------------------------------------------------------
# 2) Extend the "IEEE ICDM" Program Committees paragraph with "'18".
# ---------------------------------------------------------------------
$icdmIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p2 = $d.Paragraphs.Item($i)
    $t2 = $p2.Range.Text
    if ($t2 -like "IEEE ICDM*16*" -and $t2.Length -lt 20) {
        $icdmIndex = $i
        break
    }
}
if ($icdmIndex -eq -1) {
    throw "Could not locate the IEEE ICDM Program Committees paragraph"
}

$p2 = $d.Paragraphs.Item($icdmIndex)
$rng2 = $p2.Range
$rng2.Collapse(1)
$frag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="normal0"/><w:spacing w:line="360" w:lineRule="auto"/></w:pPr><w:r><w:t>IEEE ICDM</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>‘18</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>‘16</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng2.InsertXML($frag2)

# ---------------------------------------------------------------------
# 3) Update the cached footer PAGE field result from "7" to "2".
# ---------------------------------------------------------------------
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)
    $ftr = $sec.Footers.Item(1)
    $ftrRng = $ftr.Range
    for ($fi = 1; $fi -le $ftrRng.Fields.Count; $fi++) {
        $fld = $ftrRng.Fields.Item($fi)
        if ($fld.Code.Text -like "*PAGE*") {
            $fld.Result.Text = "2"
        }
    }
}

Write-Output "edit complete"
